$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 102, shifting rows 102:123 down to 103:124
$ws.Rows.Item(102).Insert()

# Fill in the new record at row 102
$ws.Cells.Item(102, 1).Value = 8
$ws.Cells.Item(102, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(102, 3).Value = "Coquimbo"
$ws.Cells.Item(102, 4).Value = 44711
$ws.Cells.Item(102, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(102, 5).Value = 4
$ws.Cells.Item(102, 6).Value = "Fruta"
$ws.Cells.Item(102, 7).Value = 100109
$ws.Cells.Item(102, 8).Value = "Uva"
$ws.Cells.Item(102, 9).Value = 100109001
$ws.Cells.Item(102, 10).Value = "Uva"
$ws.Cells.Item(102, 11).Value = "Red Globe"
$ws.Cells.Item(102, 12).Value = "Primera"
$ws.Cells.Item(102, 13).Value = 400
$ws.Cells.Item(102, 14).Value = 8000
$ws.Cells.Item(102, 15).Value = 9000
$ws.Cells.Item(102, 16).Value = 8500
$ws.Cells.Item(102, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(102, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(102, 19).Value = 472
$ws.Cells.Item(102, 20).Value = 18
